# 项目计划表.xlsx — "Add files via upload" edit
#
# Summary of the change being reproduced (per the OOXML diff):
#   - Column C (完成情况/"completion status") for the five plan rows (rows 3-7)
#     gets a value of 1, formatted as a percentage (-> displays "100%").
#     This introduces one new cell style (numFmtId 9 = built-in "0%").
#   - D5 (备注/"remarks" for 林玮成's row) gets new note text.
#   - A8's "总结：" ("Summary:") label is re-entered.
#   - Column D is widened (24 -> ~31.5 chars) to fit the new remark text.
#   - The active selection left by the editing session is F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3-7: mark "完成情况" (completion) as 100% -----------------------
$doneRange = $ws.Range("C3:C7")
$doneRange.Value = 1
$doneRange.NumberFormat = "0%"

# --- D5: add the remark that the use-case diagram was merged -------------
$ws.Range("D5").Value = "将用例图归并到管理群用例图"

# --- A8: (re)enter the "总结：" summary label -----------------------------
$ws.Range("A8").Value = "总结："

# --- Column D: widen to fit the longer remark text ------------------------
# The saved `width` column attribute is derived from `ColumnWidth` through
# the host's pixel grid: width = (Round(ColumnWidth * 7) + 5) / 7 (7 = the
# default font's max digit width). To land the serialized width on 31.5
# characters we invert that formula (31.5 - 5/7).
$ws.Columns.Item(4).ColumnWidth = 31.5 - 5 / 7

# --- Leave the selection where the author's session ended it -------------
$ws.Range("F4").Select()
